# Updates the "cryptos" price table (columns D = Price, E = Volume(1h))
# to the refreshed values from the latest GitHub Actions scrape.
# Leading "'" on numeric-looking D values forces Excel to store them
# as text (matching the original inlineStr cell type) instead of
# auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '30.542.60'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '2.109.29'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("D4").Value = '''1.011'
$ws.Range("D5").Value = '''335.77'
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("D7").Value = '''0.5231'
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").Value = '''0.4538'
$ws.Range("E8").Value = '  +3.99%  '
$ws.Range("D9").Value = '''55.34'
$ws.Range("E9").Value = '  +2.40%  '
$ws.Range("D10").Value = '''0.09072'
$ws.Range("E10").Value = '  +2.72%  '
$ws.Range("D11").Value = '''1.169'
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("D12").Value = '''24.53'
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("D13").Value = '2.110.06'
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("D14").Value = '''6.829'
$ws.Range("E14").Value = '  +2.35%  '
$ws.Range("D15").Value = '''8.091'
$ws.Range("E15").Value = '  +5.67%  '
$ws.Range("D16").Value = '''0.00001169'
$ws.Range("E16").Value = '  +4.60%  '
$ws.Range("D17").Value = '''96.94'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").Value = '''0.06681'
$ws.Range("E19").Value = '  +1.35%  '
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D22").Value = '''6.254'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '30.586.22'
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").Value = '''12.78'
$ws.Range("E24").Value = '  +4.65%  '
$ws.Range("D25").Value = '''2.355'
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("D26").Value = '2.353.78'
$ws.Range("E26").Value = '  +0.86%  '
$ws.Range("D27").Value = '''22.23'
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").Value = '''163.71'
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("D29").Value = '''2.512'
$ws.Range("E29").Value = '  -1.22%  '
$ws.Range("D30").Value = '''133.41'
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("E31").Value = '  +2.55%  '
$ws.Range("D32").Value = '''0.1067'
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").Value = '''1.639'
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("D34").Value = '''6.345'
$ws.Range("E34").Value = '  +3.40%  '
$ws.Range("D35").Value = '''3.959'
$ws.Range("E35").Value = '  +1.39%  '
$ws.Range("D36").Value = '''10.42'
$ws.Range("E36").Value = '  +2.98%  '
$ws.Range("D37").Value = '''5.900'
$ws.Range("D38").Value = '''0.02613'
$ws.Range("E38").Value = '  +1.65%  '
$ws.Range("D39").Value = '''0.06806'
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").Value = '''0.2313'
$ws.Range("E40").Value = '  +2.96%  '
$ws.Range("D41").Value = '''12.59'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '''0.6856'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = '''1.257'
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D44").Value = '''14.66'
$ws.Range("E44").Value = '  +5.19%  '
$ws.Range("E45").Value = '  +1.61%  '
$ws.Range("D46").Value = '''2.295'
$ws.Range("E46").Value = '  +4.90%  '
$ws.Range("D47").Value = '''3.685'
$ws.Range("E47").Value = '  +1.59%  '
$ws.Range("D48").Value = '''0.00000000353'
$ws.Range("E48").Value = '  +18.69%  '
$ws.Range("D49").Value = '''1.251'
$ws.Range("E49").Value = '  +0.99%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '''82.99'
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = '''0.3373'
$ws.Range("E51").Value = '  +13.10%  '
